$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Wealth Row: 6 -> 7
$ws.Range("B5").Value = 7

# Records Row: 2 -> 3
$ws.Range("B6").Value = 3

# Records Banks Column: "G" -> "J"
$ws.Range("B7").Value = "J"
